$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D37").Value = "[Paper Review] Latent Space Autoregression for Novelty Detection"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1438&mod=document&pageid=1"

$ws.Range("D42").Value = "Python 인스타그램 크롤링 환경 셋팅"
$ws.Range("E42").Value = "https://kjk92.tistory.com/62"

$ws.Range("D44").Value = "Hyper-Parameter Tuning 및 AutoML 논문 리뷰"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/75"
